$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 91.63362566666667
$ws.Cells.Item(2, 8).Value = 274.900877
$ws.Cells.Item(2, 9).Value = 0.5385978585809309
$ws.Cells.Item(2, 10).Value = 0.538597858580931
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 6.490547666666667
$ws.Cells.Item(2, 14).Value = 19.471643
$ws.Cells.Item(2, 15).Value = 0.8021666724616637
$ws.Cells.Item(2, 16).Value = 0.8021666724616636
$ws.Cells.Item(2, 17).Value = 594.7524152589901
$ws.Cells.Item(2, 18).Value = 5352.771737330911
$ws.Cells.Item(2, 19).Value = 0.4320452520128431
$ws.Cells.Item(2, 20).Value = 0.4320452520128431

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 91.63362566666667
$ws.Cells.Item(3, 8).Value = 274.900877
$ws.Cells.Item(3, 9).Value = 0.5385978585809309
$ws.Cells.Item(3, 10).Value = 0.538597858580931
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 1.600723
$ws.Cells.Item(3, 14).Value = 4.802169
$ws.Cells.Item(3, 15).Value = 0.1978333275383364
$ws.Cells.Item(3, 16).Value = 0.1978333275383364
$ws.Cells.Item(3, 17).Value = 146.6800521780237
$ws.Cells.Item(3, 18).Value = 1320.120469602213
$ws.Cells.Item(3, 19).Value = 0.1065526065680879
$ws.Cells.Item(3, 20).Value = 0.1065526065680879

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 2.089075
$ws.Cells.Item(4, 8).Value = 6.267225
$ws.Cells.Item(4, 9).Value = 0.01227902217367198
$ws.Cells.Item(4, 10).Value = 0.01227902217367199
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 6.490547666666667
$ws.Cells.Item(4, 14).Value = 19.471643
$ws.Cells.Item(4, 15).Value = 0.8021666724616637
$ws.Cells.Item(4, 16).Value = 0.8021666724616636
$ws.Cells.Item(4, 17).Value = 13.55924086674167
$ws.Cells.Item(4, 18).Value = 122.033167800675
$ws.Cells.Item(4, 19).Value = 0.009849822358137439
$ws.Cells.Item(4, 20).Value = 0.009849822358137441

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 2.089075
$ws.Cells.Item(5, 8).Value = 6.267225
$ws.Cells.Item(5, 9).Value = 0.01227902217367198
$ws.Cells.Item(5, 10).Value = 0.01227902217367199
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 1.600723
$ws.Cells.Item(5, 14).Value = 4.802169
$ws.Cells.Item(5, 15).Value = 0.1978333275383364
$ws.Cells.Item(5, 16).Value = 0.1978333275383364
$ws.Cells.Item(5, 17).Value = 3.344030401225
$ws.Cells.Item(5, 18).Value = 30.096273611025
$ws.Cells.Item(5, 19).Value = 0.002429199815534545
$ws.Cells.Item(5, 20).Value = 0.002429199815534545

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 74.500838
$ws.Cells.Item(6, 8).Value = 223.502514
$ws.Cells.Item(6, 9).Value = 0.437895930858942
$ws.Cells.Item(6, 10).Value = 0.4378959308589421
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 6.490547666666667
$ws.Cells.Item(6, 14).Value = 19.471643
$ws.Cells.Item(6, 15).Value = 0.8021666724616637
$ws.Cells.Item(6, 16).Value = 0.8021666724616636
$ws.Cells.Item(6, 17).Value = 483.5512402456114
$ws.Cells.Item(6, 18).Value = 4351.961162210502
$ws.Cells.Item(6, 19).Value = 0.3512655217416203
$ws.Cells.Item(6, 20).Value = 0.3512655217416203

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 74.500838
$ws.Cells.Item(7, 8).Value = 223.502514
$ws.Cells.Item(7, 9).Value = 0.437895930858942
$ws.Cells.Item(7, 10).Value = 0.4378959308589421
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 1.600723
$ws.Cells.Item(7, 14).Value = 4.802169
$ws.Cells.Item(7, 15).Value = 0.1978333275383364
$ws.Cells.Item(7, 16).Value = 0.1978333275383364
$ws.Cells.Item(7, 17).Value = 119.255204905874
$ws.Cells.Item(7, 18).Value = 1073.296844152866
$ws.Cells.Item(7, 19).Value = 0.08663040911732178
$ws.Cells.Item(7, 20).Value = 0.08663040911732178

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 1.910122666666666
$ws.Cells.Item(8, 8).Value = 5.730367999999999
$ws.Cells.Item(8, 9).Value = 0.01122718838645499
$ws.Cells.Item(8, 10).Value = 0.01122718838645499
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 6.490547666666667
$ws.Cells.Item(8, 14).Value = 19.471643
$ws.Cells.Item(8, 15).Value = 0.8021666724616637
$ws.Cells.Item(8, 16).Value = 0.8021666724616636
$ws.Cells.Item(8, 17).Value = 12.39774221718044
$ws.Cells.Item(8, 18).Value = 111.579679954624
$ws.Cells.Item(8, 19).Value = 0.009006076349062834
$ws.Cells.Item(8, 20).Value = 0.009006076349062834

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 1.910122666666666
$ws.Cells.Item(9, 8).Value = 5.730367999999999
$ws.Cells.Item(9, 9).Value = 0.01122718838645499
$ws.Cells.Item(9, 10).Value = 0.01122718838645499
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 1.600723
$ws.Cells.Item(9, 14).Value = 4.802169
$ws.Cells.Item(9, 15).Value = 0.1978333275383364
$ws.Cells.Item(9, 16).Value = 0.1978333275383364
$ws.Cells.Item(9, 17).Value = 3.057577285354666
$ws.Cells.Item(9, 18).Value = 27.518195568192
$ws.Cells.Item(9, 19).Value = 0.002221112037392157
$ws.Cells.Item(9, 20).Value = 0.002221112037392157
